$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Range("H38").Value = 3641.238
$ws.Range("J38").Value = 5749.0835
$ws.Range("L38").Value = 17247.2505
$ws.Range("N38").Value = -17991.2505
# Row 40
$ws.Range("H40").Value = 2475.7856
$ws.Range("I40").Value = 1159.3334
$ws.Range("K40").Value = 1159.3334
$ws.Range("M40").Value = -984.3334
# Row 62
$ws.Range("H62").Value = 53555.234
$ws.Range("I62").Value = 70869.914
$ws.Range("K62").Value = 70869.914
$ws.Range("M62").Value = -70245.914
# Row 65
$ws.Range("H65").Value = 53555.234
$ws.Range("I65").Value = 70869.914
$ws.Range("K65").Value = 354349.57
$ws.Range("M65").Value = -351229.57
# Row 112
$ws.Range("H112").Value = 2349.3
$ws.Range("I112").Value = 1674.6666
$ws.Range("J112").Value = 2638.4285
$ws.Range("K112").Value = 5023.9998
$ws.Range("L112").Value = 7915.2855
$ws.Range("M112").Value = -3915.9998
$ws.Range("N112").Value = -10131.2855
# Row 117
$ws.Range("H117").Value = 75000
$ws.Range("J117").Value = 75000
$ws.Range("L117").Value = 75000
$ws.Range("N117").Value = -84178
# Row 127
$ws.Range("H127").Value = 1433.8
$ws.Range("I127").Value = 1410
$ws.Range("J127").Value = 1489.3334
$ws.Range("K127").Value = 4230
$ws.Range("L127").Value = 4468.0002
$ws.Range("M127").Value = 730
$ws.Range("N127").Value = -14388.0002
# Row 132
$ws.Range("H132").Value = 2152.6667
$ws.Range("I132").Value = 1545.8334
$ws.Range("J132").Value = 3366.3333
$ws.Range("K132").Value = 4637.5002
$ws.Range("L132").Value = 10098.9999
$ws.Range("M132").Value = -2107.5002
$ws.Range("N132").Value = -15158.9999
# Row 137
$ws.Range("H137").Value = 1371.2307
$ws.Range("I137").Value = 1183.5
$ws.Range("K137").Value = 3550.5
$ws.Range("M137").Value = -1000.5
# Row 138
$ws.Range("H138").Value = 3067.5833
$ws.Range("I138").Value = 2212.6155
$ws.Range("J138").Value = 3304.0637
$ws.Range("K138").Value = 6637.8465
$ws.Range("L138").Value = 9912.1911
$ws.Range("M138").Value = -1497.8465
$ws.Range("N138").Value = -20192.1911

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5195.115
$ws.Range("I32").Value = 5197.619
$ws.Range("J32").Value = 5184.6
$ws.Range("K32").Value = 5197.619
$ws.Range("L32").Value = 5184.6
$ws.Range("M32").Value = -4910.619
$ws.Range("N32").Value = -5758.6
# Row 61
$ws.Range("H61").Value = 6863.091
$ws.Range("I61").Value = 6863.091
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 6863.091
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -6651.091
$ws.Range("N61").Value = ""
# Row 74
$ws.Range("H74").Value = 1261.7097
$ws.Range("I74").Value = 1261.7097
$ws.Range("K74").Value = 1261.7097
$ws.Range("M74").Value = -387.7097000000001
# Row 77
$ws.Range("H77").Value = 1261.7097
$ws.Range("I77").Value = 1261.7097
$ws.Range("K77").Value = 6308.548500000001
$ws.Range("M77").Value = -1940.548500000001
# Row 97
$ws.Range("H97").Value = 499.7647
$ws.Range("I97").Value = 302.125
$ws.Range("J97").Value = 675.44446
$ws.Range("K97").Value = 302.125
$ws.Range("L97").Value = 675.44446
$ws.Range("M97").Value = 193.875
$ws.Range("N97").Value = -1667.44446
# Row 136
$ws.Range("H136").Value = 6863.091
$ws.Range("I136").Value = 6863.091
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 20589.273
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -18039.273
$ws.Range("N136").Value = ""

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 3133.182
$ws.Range("J58").Value = 2243.8333
$ws.Range("L58").Value = 2243.8333
$ws.Range("N58").Value = -2649.8333
# Row 136
$ws.Range("H136").Value = 3133.182
$ws.Range("J136").Value = 2243.8333
$ws.Range("L136").Value = 6731.499899999999
$ws.Range("N136").Value = -11831.4999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 16
$ws.Range("H16").Value = 1000
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").Value = ""
# Row 33
$ws.Range("H33").Value = 148.18182
$ws.Range("J33").Value = 224.8
$ws.Range("L33").Value = 1348.8
$ws.Range("N33").Value = -1914.8
# Row 34
$ws.Range("H34").Value = 3762.1082
$ws.Range("J34").Value = 3970.5293
$ws.Range("L34").Value = 11911.5879
$ws.Range("N34").Value = -12079.5879
# Row 39
$ws.Range("H39").Value = 654.8
$ws.Range("I39").Value = 654.8
$ws.Range("K39").Value = 1964.4
$ws.Range("M39").Value = -1670.4
# Row 46
$ws.Range("H46").Value = 1182720.8
$ws.Range("I46").Value = 10099750
$ws.Range("J46").Value = 5000
$ws.Range("K46").Value = 30299250
$ws.Range("L46").Value = 15000
$ws.Range("M46").Value = -30299159
$ws.Range("N46").Value = -15182
# Row 49
$ws.Range("H49").Value = 1003
$ws.Range("I49").Value = 1003
$ws.Range("K49").Value = 3009
$ws.Range("M49").Value = -2853
# Row 55
$ws.Range("H55").Value = 300147.34
$ws.Range("I55").Value = 2502000
$ws.Range("J55").Value = 6567
$ws.Range("K55").Value = 7506000
$ws.Range("L55").Value = 19701
$ws.Range("M55").Value = -7505823
$ws.Range("N55").Value = -20055
# Row 64
$ws.Range("H64").Value = 1000
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = ""
# Row 67
$ws.Range("H67").Value = 1000
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = ""
# Row 70
$ws.Range("H70").Value = 541.5
$ws.Range("I70").Value = 541.5
$ws.Range("K70").Value = 1624.5
$ws.Range("M70").Value = -1309.5
# Row 73
$ws.Range("H73").Value = 541.5
$ws.Range("I73").Value = 541.5
$ws.Range("K73").Value = 1624.5
$ws.Range("M73").Value = -532.5
# Row 98
$ws.Range("H98").Value = 132.33333
$ws.Range("J98").Value = 132.33333
$ws.Range("L98").Value = 396.99999
$ws.Range("N98").Value = -3392.99999
# Row 107
$ws.Range("H107").Value = 2024.875
$ws.Range("J107").Value = 990.7273
$ws.Range("L107").Value = 2972.1819
$ws.Range("N107").Value = -6812.1819

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 4808.2607
$ws.Range("I80").Value = 3460
$ws.Range("J80").Value = 5845.385
$ws.Range("K80").Value = 3460
$ws.Range("L80").Value = 5845.385
$ws.Range("M80").Value = -2462
$ws.Range("N80").Value = -7841.385
# Row 83
$ws.Range("H83").Value = 4808.2607
$ws.Range("I83").Value = 3460
$ws.Range("J83").Value = 5845.385
$ws.Range("K83").Value = 17300
$ws.Range("L83").Value = 29226.925
$ws.Range("M83").Value = -12308
$ws.Range("N83").Value = -39210.925
# Row 119
$ws.Range("H119").Value = 72210
$ws.Range("J119").Value = 72210
$ws.Range("L119").Value = 72210
$ws.Range("N119").Value = -81886
# Row 122
$ws.Range("H122").Value = 3152.8333
$ws.Range("I122").Value = 2351.6
$ws.Range("J122").Value = 4154.375
$ws.Range("K122").Value = 7054.799999999999
$ws.Range("L122").Value = 12463.125
$ws.Range("M122").Value = -4604.799999999999
$ws.Range("N122").Value = -17363.125

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 2035.5
$ws.Range("I16").Value = 547.3333
$ws.Range("K16").Value = 547.3333
$ws.Range("M16").Value = -377.3333
# Row 40
$ws.Range("H40").Value = 4422.033
$ws.Range("I40").Value = 2876.8948
$ws.Range("J40").Value = 7090.909
$ws.Range("K40").Value = 2876.8948
$ws.Range("L40").Value = 7090.909
$ws.Range("M40").Value = -2740.8948
$ws.Range("N40").Value = -7362.909
# Row 46
$ws.Range("H46").Value = 2439.6
$ws.Range("I46").Value = 1400
$ws.Range("K46").Value = 1400
$ws.Range("M46").Value = -1212
# Row 61
$ws.Range("H61").Value = 4531.6816
$ws.Range("I61").Value = 3313.2
$ws.Range("J61").Value = 7142.7144
$ws.Range("K61").Value = 3313.2
$ws.Range("L61").Value = 7142.7144
$ws.Range("M61").Value = -3111.2
$ws.Range("N61").Value = -7546.7144
# Row 100
$ws.Range("H100").Value = 6334.6924
$ws.Range("I100").Value = 3478.7144
$ws.Range("J100").Value = 9666.666999999999
$ws.Range("K100").Value = 3478.7144
$ws.Range("L100").Value = 9666.666999999999
$ws.Range("M100").Value = -2937.7144
$ws.Range("N100").Value = -10748.667
# Row 113
$ws.Range("H113").Value = 4531.6816
$ws.Range("I113").Value = 3313.2
$ws.Range("J113").Value = 7142.7144
$ws.Range("K113").Value = 3313.2
$ws.Range("L113").Value = 7142.7144
$ws.Range("M113").Value = -1143.2
$ws.Range("N113").Value = -11482.7144
# Row 132
$ws.Range("H132").Value = 2732.5217
$ws.Range("I132").Value = 2453.5
$ws.Range("K132").Value = 7360.5
$ws.Range("M132").Value = -4830.5
# Row 136
$ws.Range("H136").Value = 23089.96
$ws.Range("I136").Value = 2154.889
$ws.Range("J136").Value = 34865.938
$ws.Range("K136").Value = 6464.667
$ws.Range("L136").Value = 104597.814
$ws.Range("M136").Value = -3914.667
$ws.Range("N136").Value = -109697.814

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 39
$ws.Range("H39").Value = 24500
$ws.Range("I39").Value = 24000
$ws.Range("K39").Value = 24000
$ws.Range("M39").Value = -23587
# Row 45
$ws.Range("H45").Value = 10500
$ws.Range("I45").Value = 11500
$ws.Range("K45").Value = 11500
$ws.Range("M45").Value = -11009
# Row 62
$ws.Range("H62").Value = 15611.885
$ws.Range("J62").Value = 16208.792
$ws.Range("L62").Value = 16208.792
$ws.Range("N62").Value = -17456.792
# Row 65
$ws.Range("H65").Value = 15611.885
$ws.Range("J65").Value = 16208.792
$ws.Range("L65").Value = 81043.95999999999
$ws.Range("N65").Value = -87283.95999999999
# Row 136
$ws.Range("H136").Value = 1627.1951
$ws.Range("I136").Value = 1114.6296
$ws.Range("J136").Value = 2615.7144
$ws.Range("K136").Value = 3343.8888
$ws.Range("L136").Value = 7847.1432
$ws.Range("M136").Value = -793.8887999999997
$ws.Range("N136").Value = -12947.1432
